$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows to append (A: Rv_ID, B: annot_int, C: Gene names, D: Function [CC], E: cluster_labels)
$data = @(
    @("Rv1030",  3, "kdpB Rv1030 MTCY10G2.19c",       "FUNCTION: Part of the high-affinity ATP-driven potassium transport (or Kdp) system, which catalyzes the hydrolysis of ATP coupled with the electrogenic transport of potassium into the cytoplasm. This subunit is responsible for energy coupling to the transport system. {ECO:0000255|HAMAP-Rule:MF_00285}.", 55),
    @("Rv1283c", 2, "Rv1283c MTCY373.02c",            "", 55),
    @("Rv0315",  1, "Rv0315",                         "", 55),
    @("Rv2836c", 1, "dinF Rv2836c",                   "", 55),
    @("Rv3218",  1, "Rv3218",                         "", 55),
    @("Rv2799",  1, "Rv2799",                         "", 55),
    @("Rv2622",  1, "Rv2622",                         "", 55),
    @("Rv1211",  1, "Rv1211",                         "", 55),
    @("Rv0835",  1, "lpqQ Rv0835",                    "", 55)
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
